# Update MACRO_SCORE column (N) values for rows 2-5 to the refreshed figure.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = 85.92500513438651

$ws.Range("N2").Value = $newValue
$ws.Range("N3").Value = $newValue
$ws.Range("N4").Value = $newValue
$ws.Range("N5").Value = $newValue
